$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-07 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-08 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("280÷8=35, 0", $true, $false, $false, $false, $false, $true, 1, $false, "921÷4=230, 1", 2) | Out-Null
$d.Content.Find.Execute("676÷3=225, 1", $true, $false, $false, $false, $false, $true, 1, $false, "701÷3=233, 2", 2) | Out-Null
$d.Content.Find.Execute("483÷9=53, 6", $true, $false, $false, $false, $false, $true, 1, $false, "975÷6=162, 3", 2) | Out-Null
$d.Content.Find.Execute("685÷3=228, 1", $true, $false, $false, $false, $false, $true, 1, $false, "316÷2=158, 0", 2) | Out-Null
$d.Content.Find.Execute("932÷3=310, 2", $true, $false, $false, $false, $false, $true, 1, $false, "234÷6=39, 0", 2) | Out-Null
$d.Content.Find.Execute("798÷5=159, 3", $true, $false, $false, $false, $false, $true, 1, $false, "355÷7=50, 5", 2) | Out-Null
$d.Content.Find.Execute("111÷8=13, 7", $true, $false, $false, $false, $false, $true, 1, $false, "143÷7=20, 3", 2) | Out-Null
$d.Content.Find.Execute("164÷6=27, 2", $true, $false, $false, $false, $false, $true, 1, $false, "493÷2=246, 1", 2) | Out-Null
$d.Content.Find.Execute("369÷7=52, 5", $true, $false, $false, $false, $false, $true, 1, $false, "252÷9=28, 0", 2) | Out-Null
$d.Content.Find.Execute("848÷5=169, 3", $true, $false, $false, $false, $false, $true, 1, $false, "912÷2=456, 0", 2) | Out-Null
$d.Content.Find.Execute("911÷8=113, 7", $true, $false, $false, $false, $false, $true, 1, $false, "970÷6=161, 4", 2) | Out-Null
$d.Content.Find.Execute("609÷9=67, 6", $true, $false, $false, $false, $false, $true, 1, $false, "684÷7=97, 5", 2) | Out-Null
$d.Content.Find.Execute("204÷5=40, 4", $true, $false, $false, $false, $false, $true, 1, $false, "893÷4=223, 1", 2) | Out-Null
$d.Content.Find.Execute("179÷4=44, 3", $true, $false, $false, $false, $false, $true, 1, $false, "687÷6=114, 3", 2) | Out-Null
$d.Content.Find.Execute("950÷9=105, 5", $true, $false, $false, $false, $false, $true, 1, $false, "406÷7=58, 0", 2) | Out-Null
$d.Content.Find.Execute("972÷8=121, 4", $true, $false, $false, $false, $false, $true, 1, $false, "515÷7=73, 4", 2) | Out-Null
$d.Content.Find.Execute("617÷6=102, 5", $true, $false, $false, $false, $false, $true, 1, $false, "704÷2=352, 0", 2) | Out-Null
$d.Content.Find.Execute("437÷6=72, 5", $true, $false, $false, $false, $false, $true, 1, $false, "907÷3=302, 1", 2) | Out-Null
$d.Content.Find.Execute("519÷5=103, 4", $true, $false, $false, $false, $false, $true, 1, $false, "892÷7=127, 3", 2) | Out-Null
$d.Content.Find.Execute("536÷9=59, 5", $true, $false, $false, $false, $false, $true, 1, $false, "254÷5=50, 4", 2) | Out-Null
$d.Content.Find.Execute("345÷8=43, 1", $true, $false, $false, $false, $false, $true, 1, $false, "900÷7=128, 4", 2) | Out-Null
$d.Content.Find.Execute("226÷7=32, 2", $true, $false, $false, $false, $false, $true, 1, $false, "278÷8=34, 6", 2) | Out-Null
$d.Content.Find.Execute("690÷4=172, 2", $true, $false, $false, $false, $false, $true, 1, $false, "563÷9=62, 5", 2) | Out-Null
$d.Content.Find.Execute("632÷6=105, 2", $true, $false, $false, $false, $false, $true, 1, $false, "978÷3=326, 0", 2) | Out-Null
$d.Content.Find.Execute("485÷7=69, 2", $true, $false, $false, $false, $false, $true, 1, $false, "286÷8=35, 6", 2) | Out-Null
